$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A2 -> FALSE (boolean), B2 -> 137, C2 -> 17.04
$ws.Range("A2").Value = $false
$ws.Range("B2").Value = 137
$ws.Range("C2").Value = 17.04

# Row 3: A3 -> TRUE (boolean), B3 -> 667, C3 -> 82.96
$ws.Range("A3").Value = $true
$ws.Range("B3").Value = 667
$ws.Range("C3").Value = 82.95999999999999
